$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '98.082.37'
$ws.Range("E2").Value = '  +0.06%  '

$ws.Range("D3").Value = '3.384.58'
$ws.Range("E3").Value = '  -0.59%  '

$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '253.57'
$ws.Range("E5").Value = '  -0.38%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '661.63'
$ws.Range("E6").Value = '  +0.05%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.46'
$ws.Range("E7").Value = '  +1.54%  '

$ws.Range("E8").Value = '  -1.46%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.05'
$ws.Range("E9").Value = '  +0.10%  '

$ws.Range("E10").Value = '  -0.02%  '

$ws.Range("D11").Value = '3.381.77'
$ws.Range("E11").Value = '  -0.62%  '

$ws.Range("E12").Value = '  -2.45%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '41.70'
$ws.Range("E13").Value = '  -1.37%  '

$ws.Range("D14").Value = '97.833.94'
$ws.Range("E14").Value = '  +0.11%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.12'
$ws.Range("E15").Value = '  -4.99%  '

$ws.Range("E16").Value = '  -3.34%  '

$ws.Range("D17").Value = '4.014.56'
$ws.Range("E17").Value = '  -0.56%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.94'
$ws.Range("E18").Value = '  -0.27%  '

$ws.Range("D19").Value = '3.376.79'
$ws.Range("E19").Value = '  -0.78%  '

$ws.Range("E20").Value = '  +2.42%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.530'
$ws.Range("E21").Value = '  -3.93%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.95'
$ws.Range("E22").Value = '  +0.11%  '

$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '512.40'
$ws.Range("E23").Value = '  +0.90%  '

$ws.Range("B24").Value = 'SuiNetwork'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.44'
$ws.Range("E24").Value = '  +0.45%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.96'
$ws.Range("E25").Value = '  +6.52%  '

$ws.Range("E26").Value = '  -1.70%  '

$ws.Range("E27").Value = '  -3.29%  '

$ws.Range("E28").Value = '  -3.52%  '

$ws.Range("D29").Value = '3.568.17'
$ws.Range("E29").Value = '  -0.70%  '

$ws.Range("E30").Value = '  -0.91%  '

$ws.Range("E31").Value = '  -2.70%  '

$ws.Range("E32").Value = '  +0.24%  '

$ws.Range("E33").Value = '  -5.47%  '

$ws.Range("E34").Value = '  +9.43%  '

$ws.Range("E35").Value = '  +0.12%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.561'
$ws.Range("E36").Value = '  -1.93%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '28.98'
$ws.Range("E37").Value = '  -2.33%  '

$ws.Range("E38").Value = '  +0.81%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.50'
$ws.Range("E39").Value = '  +0.18%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '536.41'
$ws.Range("E40").Value = '  +0.86%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.153'
$ws.Range("E41").Value = '  +0.23%  '

$ws.Range("E42").Value = '  -0.03%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '24.40'
$ws.Range("E43").Value = '  -1.14%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.854'
$ws.Range("E44").Value = '  -2.22%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0431'
$ws.Range("E45").Value = '  +1.71%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.73'
$ws.Range("E46").Value = '  +0.30%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.69'
$ws.Range("E47").Value = '  +0.97%  '

$ws.Range("E48").Value = '  +7.71%  '

$ws.Range("B49").Value = 'OKB'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '56.18'
$ws.Range("E49").Value = '  +3.07%  '

$ws.Range("B50").Value = 'Filecoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.62'
$ws.Range("E50").Value = '  -3.47%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.60'
$ws.Range("E51").Value = '  -5.33%  '
